$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.913.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "'3.094.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'544.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "'140.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'3.087.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E13").Value = "  +3.80%  "
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'3.594.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "'63.985.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'3.095.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "'479.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "'7.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").Value = "'78.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "'12.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "'8.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.67%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'26.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'1.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "'57.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("E34").Value = "  -6.57%  "
$ws.Range("D35").Value = "'5.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.36%  "
$ws.Range("D36").Value = "'496.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.88%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'3.243.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.40%  "
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "'2.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D46").Value = "'124.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("E49").Value = "  +6.00%  "
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  +0.53%  "
